$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.637.79"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "2.473.07"
$ws.Range("E3").Value = "  -0.30%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'318.04"
$ws.Range("E5").Value = "  +1.32%  "
$ws.Range("D6").Value = "'92.63"
$ws.Range("D7").Value = "'0.553"
$ws.Range("E7").Value = "  +1.73%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'0.519"
$ws.Range("E9").Value = "  +2.60%  "
$ws.Range("D10").Value = "'0.0870"
$ws.Range("E10").Value = "  +10.97%  "
$ws.Range("D11").Value = "'32.92"
$ws.Range("E11").Value = "  +0.59%  "
$ws.Range("D12").Value = "'0.112"
$ws.Range("E12").Value = "  +1.08%  "
$ws.Range("D13").Value = "2.853.57"
$ws.Range("E13").Value = "  -0.30%  "
$ws.Range("D14").Value = "'6.90"
$ws.Range("E14").Value = "  +0.98%  "
$ws.Range("D15").Value = "'15.63"
$ws.Range("E15").Value = "  -3.07%  "
$ws.Range("D16").Value = "2.483.21"
$ws.Range("E16").Value = "  +0.65%  "
$ws.Range("D17").Value = "'0.788"
$ws.Range("E17").Value = "  +2.88%  "
$ws.Range("D18").Value = "41.597.16"
$ws.Range("E18").Value = "  +0.17%  "
$ws.Range("D19").Value = "0.0₃0956"
$ws.Range("E19").Value = "  +1.95%  "
$ws.Range("D20").Value = "'6.48"
$ws.Range("E20").Value = "  +1.26%  "
$ws.Range("D21").Value = "'71.25"
$ws.Range("E21").Value = "  -0.45%  "
$ws.Range("D22").Value = "'11.44"
$ws.Range("E22").Value = "  +1.16%  "
$ws.Range("D23").Value = "'241.06"
$ws.Range("E23").Value = "  +1.77%  "
$ws.Range("D24").Value = "'2.74"
$ws.Range("E24").Value = "  +1.32%  "
$ws.Range("E25").Value = "  +1.38%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").Value = "'24.73"
$ws.Range("E27").Value = "  -0.35%  "
$ws.Range("D28").Value = "'2.29"
$ws.Range("E28").Value = "  +3.80%  "
$ws.Range("D29").Value = "'9.90"
$ws.Range("E29").Value = "  +2.65%  "
$ws.Range("D30").Value = "'36.35"
$ws.Range("E30").Value = "  +0.92%  "
$ws.Range("D31").Value = "'157.93"
$ws.Range("E31").Value = "  -0.29%  "
$ws.Range("D32").Value = "'5.54"
$ws.Range("E32").Value = "  +1.37%  "
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("D34").Value = "'0.0774"
$ws.Range("E34").Value = "  +2.63%  "
$ws.Range("D35").Value = "'2.58"
$ws.Range("E35").Value = "  +0.38%  "
$ws.Range("D36").Value = "'17.41"
$ws.Range("E36").Value = "  +0.49%  "
$ws.Range("D37").Value = "'2.91"
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("E38").Value = "  +0.61%  "
$ws.Range("E39").Value = "  +1.43%  "
$ws.Range("E40").Value = "  -1.75%  "
$ws.Range("D41").Value = "'4.01"
$ws.Range("E41").Value = "  -3.31%  "
$ws.Range("E42").Value = "  +1.97%  "
$ws.Range("D43").Value = "1.985.69"
$ws.Range("E43").Value = "  +0.72%  "
$ws.Range("D44").Value = "'19.20"
$ws.Range("E44").Value = "  -2.10%  "
$ws.Range("E45").Value = "  +0.71%  "
$ws.Range("D46").Value = "'3.02"
$ws.Range("E46").Value = "  +1.84%  "
$ws.Range("D47").Value = "'9.24"
$ws.Range("E47").Value = "  +2.29%  "
$ws.Range("D48").Value = "2.711.31"
$ws.Range("E48").Value = "  -0.30%  "
$ws.Range("D49").Value = "'97.41"
$ws.Range("E49").Value = "  -0.36%  "
$ws.Range("D50").Value = "'67.43"
$ws.Range("E50").Value = "  -0.74%  "
$ws.Range("D51").Value = "'73.52"
$ws.Range("E51").Value = "  +1.20%  "
